# Slamf6-Slamf6.xlsx: refresh with new TPM-derived specificity values.
# The "ECs" sending/target cluster rows are dropped entirely (rows 3-5 of the
# original sheet, plus the row-2 "ECs" pairing becomes the remaining "FAPs"
# pairing), leaving a single FAPs -> FAPs data row with recomputed
# specificity figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three rows that referenced the "ECs" cluster on either side
# (old rows 3, 4 and 5); this also drops the now-unused "ECs" shared string
# and shifts the remaining sheet dimension down to A1:T2.
$ws.Rows("3:5").Delete()

# The single remaining data row (row 2) keeps the FAPs<->FAPs ligand/receptor
# pairing, but with its cluster label updated from "ECs" to "FAPs" and its
# expression/specificity metrics refreshed to the new TPM-based values.
$ws.Range("A2").Value = "FAPs"
$ws.Range("D2").Value = "FAPs"

$ws.Range("G2").Value = 0.002064333333333333
$ws.Range("H2").Value = 0.006193
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1

$ws.Range("M2").Value = 0.002064333333333333
$ws.Range("N2").Value = 0.006193
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1

$ws.Range("Q2").Value = 0.000004261472111111111
$ws.Range("R2").Value = 0.000038353249
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
